$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 3.38 = 12949.32 pesos", "1000 Bs = 3.37 = 12872.59 pesos")
$text = $text.Replace("12949.32 pesos = 3.36 = 968.68 Bs", "12872.59 pesos = 3.36 = 963.1 Bs")
$cell.Value = $text

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 297
$wsTasas.Range("O10").Value = 3823.16
$wsTasas.Range("N12").Value = 3836
$wsTasas.Range("O12").Value = 287
